# labelled_M_0098_09y8m_1_fa.xlsx — "almost done 2 flies left to clean"
#
# Turns several of the plain C-column "range start" values into formulas
# that reference the prior row's D-column "range end" (C[n] = D[n-1]),
# and nudges a handful of D-column boundary values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- D-column value edits (do these first so dependent C formulas below
#     recalculate against the final numbers) ---------------------------
$ws.Range("D34").Value  = 22.85
$ws.Range("D35").Value  = 23.43
$ws.Range("D70").Value  = 44.8
$ws.Range("D95").Value  = 69.254999999999995
$ws.Range("D107").Value = 77.849999999999994
$ws.Range("D166").Value = 120.075

# --- C-column values that shift without becoming formulas -------------
$ws.Range("C180").Value = 129.59
$ws.Range("D180").Value = 130.29
$ws.Range("C181").Value = 130.29
$ws.Range("D181").Value = 130.59
$ws.Range("C182").Value = 130.59

# --- C32:C41 become formulas pointing at the row above's D value -------
$ws.Range("C32").Formula = "=D31"
$ws.Range("C33").Formula = "=D32"
$ws.Range("C34").Formula = "=D33"
$ws.Range("C35").Formula = "=D34"
$ws.Range("C36").Formula = "=D35"
$ws.Range("C37").Formula = "=D36"
$ws.Range("C38").Formula = "=D37"
$ws.Range("C39").Formula = "=D38"
$ws.Range("C40").Formula = "=D39"
$ws.Range("C41").Formula = "=D40"

# --- C96:C99 become formulas the same way -------------------------------
$ws.Range("C96").Formula = "=D95"
$ws.Range("C97").Formula = "=D96"
$ws.Range("C98").Formula = "=D97"
$ws.Range("C99").Formula = "=D98"

# --- sheet view: zoom + selection (topLeftCell isn't exposed by this
#     COM host, best-effort via ScrollRow/ScrollColumn) -----------------
$excel.ActiveWindow.ScrollRow = 180
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.Zoom = 118
$null = $ws.Range("D183").Select()
